$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.5867285
$ws.Range("H2").Value = 1.173457
$ws.Range("I2").Value = 0.1699339543088995
$ws.Range("J2").Value = 0.1343654643255494
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.5867285
$ws.Range("N2").Value = 1.173457
$ws.Range("O2").Value = 0.1699339543088995
$ws.Range("P2").Value = 0.1343654643255494
$ws.Range("Q2").Value = 0.34425033271225
$ws.Range("R2").Value = 1.377001330849
$ws.Range("S2").Value = 0.02887754882705916
$ws.Range("T2").Value = 0.01805407800342048

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.5867285
$ws.Range("H3").Value = 1.173457
$ws.Range("I3").Value = 0.1699339543088995
$ws.Range("J3").Value = 0.1343654643255494
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.282443666666667
$ws.Range("N3").Value = 3.847331000000001
$ws.Range("O3").Value = 0.3714336758058812
$ws.Range("P3").Value = 0.4405346052127009
$ws.Range("Q3").Value = 0.7524462488778334
$ws.Range("R3").Value = 4.514677493267
$ws.Range("S3").Value = 0.06311919329318322
$ws.Range("T3").Value = 0.05919263678087714

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 0.5867285
$ws.Range("H4").Value = 1.173457
$ws.Range("I4").Value = 0.1699339543088995
$ws.Range("J4").Value = 0.1343654643255494
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5455083333333333
$ws.Range("N4").Value = 1.636525
$ws.Range("O4").Value = 0.1579953729736847
$ws.Range("P4").Value = 0.1873885804979387
$ws.Range("Q4").Value = 0.3200652861541666
$ws.Range("R4").Value = 1.920391716925
$ws.Range("S4").Value = 0.02684877849192768
$ws.Range("T4").Value = 0.02517855362791113

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 0.5867285
$ws.Range("H5").Value = 1.173457
$ws.Range("I5").Value = 0.1699339543088995
$ws.Range("J5").Value = 0.1343654643255494
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 1.038005
$ws.Range("N5").Value = 2.07601
$ws.Range("O5").Value = 0.3006369969115345
$ws.Range("P5").Value = 0.237711349963811
$ws.Range("Q5").Value = 0.6090271166425
$ws.Range("R5").Value = 2.43610846657
$ws.Range("S5").Value = 0.05108843369672948
$ws.Range("T5").Value = 0.03194019591334064

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.282443666666667
$ws.Range("H6").Value = 3.847331000000001
$ws.Range("I6").Value = 0.3714336758058812
$ws.Range("J6").Value = 0.4405346052127009
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 0.5867285
$ws.Range("N6").Value = 1.173457
$ws.Range("O6").Value = 0.1699339543088995
$ws.Range("P6").Value = 0.1343654643255494
$ws.Range("Q6").Value = 0.7524462488778334
$ws.Range("R6").Value = 4.514677493267
$ws.Range("S6").Value = 0.06311919329318322
$ws.Range("T6").Value = 0.05919263678087714

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.282443666666667
$ws.Range("H7").Value = 3.847331000000001
$ws.Range("I7").Value = 0.3714336758058812
$ws.Range("J7").Value = 0.4405346052127009
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.282443666666667
$ws.Range("N7").Value = 3.847331000000001
$ws.Range("O7").Value = 0.3714336758058812
$ws.Range("P7").Value = 0.4405346052127009
$ws.Range("Q7").Value = 1.644661758173445
$ws.Range("R7").Value = 14.801955823561
$ws.Range("S7").Value = 0.1379629755226684
$ws.Range("T7").Value = 0.1940707383899102

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.282443666666667
$ws.Range("H8").Value = 3.847331000000001
$ws.Range("I8").Value = 0.3714336758058812
$ws.Range("J8").Value = 0.4405346052127009
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.5455083333333333
$ws.Range("N8").Value = 1.636525
$ws.Range("O8").Value = 0.1579953729736847
$ws.Range("P8").Value = 0.1873885804979387
$ws.Range("Q8").Value = 0.6995837071972222
$ws.Range("R8").Value = 6.296253364775
$ws.Range("S8").Value = 0.0586848021439369
$ws.Range("T8").Value = 0.08255115433102785

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.282443666666667
$ws.Range("H9").Value = 3.847331000000001
$ws.Range("I9").Value = 0.3714336758058812
$ws.Range("J9").Value = 0.4405346052127009
$ws.Range("K9").Value = 2
$ws.Range("M9").Value = 1.038005
$ws.Range("N9").Value = 2.07601
$ws.Range("O9").Value = 0.3006369969115345
$ws.Range("P9").Value = 0.237711349963811
$ws.Range("Q9").Value = 1.331182938218334
$ws.Range("R9").Value = 7.987097629310002
$ws.Range("S9").Value = 0.1116667048460926
$ws.Range("T9").Value = 0.1047200757108857

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.5455083333333333
$ws.Range("H10").Value = 1.636525
$ws.Range("I10").Value = 0.1579953729736847
$ws.Range("J10").Value = 0.1873885804979387
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 0.5867285
$ws.Range("N10").Value = 1.173457
$ws.Range("O10").Value = 0.1699339543088995
$ws.Range("P10").Value = 0.1343654643255494
$ws.Range("Q10").Value = 0.3200652861541666
$ws.Range("R10").Value = 1.920391716925
$ws.Range("S10").Value = 0.02684877849192768
$ws.Range("T10").Value = 0.02517855362791113

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.5455083333333333
$ws.Range("H11").Value = 1.636525
$ws.Range("I11").Value = 0.1579953729736847
$ws.Range("J11").Value = 0.1873885804979387
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.282443666666667
$ws.Range("N11").Value = 3.847331000000001
$ws.Range("O11").Value = 0.3714336758058812
$ws.Range("P11").Value = 0.4405346052127009
$ws.Range("Q11").Value = 0.6995837071972222
$ws.Range("R11").Value = 6.296253364775
$ws.Range("S11").Value = 0.0586848021439369
$ws.Range("T11").Value = 0.08255115433102785

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.5455083333333333
$ws.Range("H12").Value = 1.636525
$ws.Range("I12").Value = 0.1579953729736847
$ws.Range("J12").Value = 0.1873885804979387
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.5455083333333333
$ws.Range("N12").Value = 1.636525
$ws.Range("O12").Value = 0.1579953729736847
$ws.Range("P12").Value = 0.1873885804979387
$ws.Range("Q12").Value = 0.2975793417361111
$ws.Range("R12").Value = 2.678214075624999
$ws.Range("S12").Value = 0.02496253788109375
$ws.Range("T12").Value = 0.03511448010103246

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.5455083333333333
$ws.Range("H13").Value = 1.636525
$ws.Range("I13").Value = 0.1579953729736847
$ws.Range("J13").Value = 0.1873885804979387
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 1.038005
$ws.Range("N13").Value = 2.07601
$ws.Range("O13").Value = 0.3006369969115345
$ws.Range("P13").Value = 0.237711349963811
$ws.Range("Q13").Value = 0.5662403775416667
$ws.Range("R13").Value = 3.39744226525
$ws.Range("S13").Value = 0.0474992544567264
$ws.Range("T13").Value = 0.04454439243796728

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 1.038005
$ws.Range("H14").Value = 2.07601
$ws.Range("I14").Value = 0.3006369969115345
$ws.Range("J14").Value = 0.237711349963811
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 0.5867285
$ws.Range("N14").Value = 1.173457
$ws.Range("O14").Value = 0.1699339543088995
$ws.Range("P14").Value = 0.1343654643255494
$ws.Range("Q14").Value = 0.6090271166425
$ws.Range("R14").Value = 2.43610846657
$ws.Range("S14").Value = 0.05108843369672948
$ws.Range("T14").Value = 0.03194019591334064

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 1.038005
$ws.Range("H15").Value = 2.07601
$ws.Range("I15").Value = 0.3006369969115345
$ws.Range("J15").Value = 0.237711349963811
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.282443666666667
$ws.Range("N15").Value = 3.847331000000001
$ws.Range("O15").Value = 0.3714336758058812
$ws.Range("P15").Value = 0.4405346052127009
$ws.Range("Q15").Value = 1.331182938218334
$ws.Range("R15").Value = 7.987097629310002
$ws.Range("S15").Value = 0.1116667048460926
$ws.Range("T15").Value = 0.1047200757108857

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 1.038005
$ws.Range("H16").Value = 2.07601
$ws.Range("I16").Value = 0.3006369969115345
$ws.Range("J16").Value = 0.237711349963811
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.5455083333333333
$ws.Range("N16").Value = 1.636525
$ws.Range("O16").Value = 0.1579953729736847
$ws.Range("P16").Value = 0.1873885804979387
$ws.Range("Q16").Value = 0.5662403775416667
$ws.Range("R16").Value = 3.39744226525
$ws.Range("S16").Value = 0.0474992544567264
$ws.Range("T16").Value = 0.04454439243796728

$ws.Range("E17").Value = 2
$ws.Range("G17").Value = 1.038005
$ws.Range("H17").Value = 2.07601
$ws.Range("I17").Value = 0.3006369969115345
$ws.Range("J17").Value = 0.237711349963811
$ws.Range("K17").Value = 2
$ws.Range("M17").Value = 1.038005
$ws.Range("N17").Value = 2.07601
$ws.Range("O17").Value = 0.3006369969115345
$ws.Range("P17").Value = 0.237711349963811
$ws.Range("Q17").Value = 1.077454380025
$ws.Range("R17").Value = 4.3098175201
$ws.Range("S17").Value = 0.09038260391198601
$ws.Range("T17").Value = 0.05650668590161743
